$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the value in C3: "sediment thickness" -> "sediment_thickness"
$ws.Range("C3").Value = "sediment_thickness"

# Update the active selection from C7 to C3
$ws.Range("C3").Select()
